# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G on Sheet1 is labeled "K" (strikeouts). This script rewrites the
# "K" values for each outing row (rows 2-45) with the regenerated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 2
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 2
    16 = 0
    17 = 0
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 2
    26 = 1
    27 = 2
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 0
    35 = 1
    36 = 1
    37 = 0
    38 = 0
    39 = 1
    40 = 1
    42 = 1
    45 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
